# Update "Mining Profits Analysis.xlsx" to reflect Jun 12 data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenario Analysis")
$ws.Activate()

# "Prepared on:" date moves from May 24 2021 (44340) to Jun 12 2021 (44359)
$ws.Range("C3").Value = 44359

# ETH Price (CAD) updated
$ws.Range("I5").Value = 2778.17

# Network Diff (TH/s) updated
$ws.Range("M5").Value = 13.21

# Reflect the cell the author had selected last (B28:Q28) instead of the old S6
$ws.Range("B28:Q28").Select()
